# Etelis snappers Migrate Hypothesis AML_EDC.xlsx
# "remove the bad apple from E. coruscans and lots of other revisions"
#
# 1. On the "carbunculus" sheet (sheet2), append a new "Pure gene flow full"
#    table (rows 48-52), mirroring the existing "Pure gene flow N" tables but
#    with every cell in the migration matrix set to "*" (full/unrestricted
#    gene flow between all three populations).
# 2. Make "carbunculus" the active tab/sheet (it was "coruscans" before),
#    and leave the new table's header row selected on that sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # coruscans
$ws2 = $wb.Worksheets.Item(2)   # carbunculus

# --- add the new "Pure gene flow full" block under the existing tables ----
$ws2.Range("A48").Value2 = "Pure gene flow full"

# Copy the JA/MHI/NWHI x JA/MHI/NWHI header + row-label layout from the
# "Pure gene flow 1" block (rows 39:42) so styles/borders match exactly.
$ws2.Range("B39:E42").Copy($ws2.Range("B49:E52")) | Out-Null

# This is the "full" gene-flow model: every population pair (including the
# diagonal) uses the custom-migration wildcard "*", unlike the other
# "Pure gene flow" variants which pin some pairs to 0.
$ws2.Range("C50:E52").Value2 = "*"

# --- view/selection state ---------------------------------------------
# carbunculus becomes the active sheet/tab, coruscans stays scrolled where
# it was left (selection untouched) but is no longer the selected tab.
$ws2.Activate() | Out-Null
$ws2.Range("G48").Select() | Out-Null
